$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.834.06'
$ws.Range('E2').Value = '  -1.31%  '
$ws.Range('D3').Value = '3.901.11'
$ws.Range('E3').Value = '  -1.95%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '606.01'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '179.38'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +6.86%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.671'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.65%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.755'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.27%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.179'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +7.20%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '54.38'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.93%  '
$ws.Range('E12').Value = '  +2.82%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.57'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.56%  '
$ws.Range('D14').Value = '4.536.21'
$ws.Range('E14').Value = '  -1.79%  '
$ws.Range('D15').Value = '3.903.70'
$ws.Range('E15').Value = '  -2.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.77'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.51%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.02'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.06%  '
$ws.Range('E18').Value = '  -4.62%  '
$ws.Range('E19').Value = '  -2.08%  '
$ws.Range('D20').Value = '71.707.78'
$ws.Range('E20').Value = '  -1.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '442.21'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.28%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.82'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.25%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '94.62'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.57%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.27'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.26%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.98'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.57%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.70'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.33%  '
$ws.Range('E27').Value = '  -6.46%  '
$ws.Range('E28').Value = '  +0.85%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.51'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.92%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.77'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +13.33%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '35.38'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.01%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '13.70'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.35%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '47.86'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.97%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.126'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.42%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0000100'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +13.52%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '69.96'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.56%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '638.95'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.432'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.25%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.147'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.90%  '
$ws.Range('B40').Value = 'Dai'
$ws.Range('C40').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.00'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.21%  '
$ws.Range('B41').Value = 'ThetaToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.33'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.91%  '
$ws.Range('E42').Value = '  -0.17%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.20'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +17.57%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0473'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.20%  '
$ws.Range('B45').Value = 'Fetch.AI'
$ws.Range('C45').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.84'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +7.06%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.30'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.14%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.144'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.40%  '
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').Value = '2.940.98'
$ws.Range('E48').Value = '  +1.02%  '
$ws.Range('B49').Value = 'FLOKI'
$ws.Range('C49').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.000282'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +6.53%  '
$ws.Range('B50').Value = 'ApeXProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.27'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.16%  '
$ws.Range('B51').Value = 'WEMIXToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.77'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -18.58%  '
